# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear C2 entirely, update E2
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 8.045645122021927

# Row 3
$ws.Range("C3").Value = -5.232639093663827

# Row 4
$ws.Range("C4").Value = 2.124540184802992

# Row 5
$ws.Range("C5").Value = 8.081020954067753

# Row 6
$ws.Range("C6").Value = 4.489210662380971

# Row 7
$ws.Range("C7").Value = -0.8752093743685241

# Row 9
$ws.Range("C9").Value = 4.818339085077561

# Row 11
$ws.Range("C11").Value = 4.613634856640769

# Row 13
$ws.Range("E13").Value = 4.838485897465628

# Row 15
$ws.Range("E15").Value = 0.869978169785246

# Row 16
$ws.Range("E16").Value = 2.551560717335266

# Row 18
$ws.Range("C18").Value = -2.447533648174649

# Row 19
$ws.Range("C19").Value = 1.038949519463617
$ws.Range("E19").Value = -1.220869074712128
